# "Generate Report for Handoff"
#
# The 08ecc5f1-...md file has finished its "Handed back" cycle and is being
# sent out for a new handoff round, so it is moved to the bottom of each
# sheet's rotation (rows 2-4) and its status flips to "Ready for handoff"
# with fresh handoff/handback timestamps. The other two files each shift up
# one row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A2").Value = "ffff2779862e-1c90-401b-9786-9e4d2cf9062d.md"
$ws1.Range("A3").Value = "ffffff4d4b21dd-2fbd-4145-b0e7-ea00cb7bc0b8.md"

$ws1.Range("A4").Value = "08ecc5f1-0974-4d29-aea9-f8af491d9dda.md"
$ws1.Range("B4").Value = "Ready for handoff"
$ws1.Range("C4").Value = "Ready for handoff"

# ---------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A2").Value = "ffff2779862e-1c90-401b-9786-9e4d2cf9062d.md"
$ws2.Range("C2").Value = "a7e92d36-29de-4fd2-b571-1cfa984205c8.aa6171980a82ab00837cc07c12fda249f142ae56.zh-cn.xlf"
$ws2.Range("D2").Value = "2016-01-28 04:28:03"
$ws2.Range("E2").Value = "a7e92d36-29de-4fd2-b571-1cfa984205c8.md"
$ws2.Range("F2").Value = "a7e92d36-29de-4fd2-b571-1cfa984205c8.aa6171980a82ab00837cc07c12fda249f142ae56.zh-cn.xlf"
$ws2.Range("G2").Value = "2016-01-28 04:28:42"

$ws2.Range("A3").Value = "ffffff4d4b21dd-2fbd-4145-b0e7-ea00cb7bc0b8.md"

$ws2.Range("A4").Value = "08ecc5f1-0974-4d29-aea9-f8af491d9dda.md"
$ws2.Range("B4").Value = "Ready for handoff"
$ws2.Range("C4").Value = "08ecc5f1-0974-4d29-aea9-f8af491d9dda.4f1445b64f0f841ffe719951dfaafccd1dce7119.zh-cn.xlf"
$ws2.Range("D4").Value = "2016-01-28 04:33:33"
$ws2.Range("E4").Value = "08ecc5f1-0974-4d29-aea9-f8af491d9dda.md"
$ws2.Range("F4").Value = "08ecc5f1-0974-4d29-aea9-f8af491d9dda.4f1445b64f0f841ffe719951dfaafccd1dce7119.zh-cn.xlf"
$ws2.Range("G4").Value = "2016-01-28 04:32:36"

# ---------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A2").Value = "ffff2779862e-1c90-401b-9786-9e4d2cf9062d.md"
$ws3.Range("C2").Value = "a7e92d36-29de-4fd2-b571-1cfa984205c8.aa6171980a82ab00837cc07c12fda249f142ae56.de-de.xlf"
$ws3.Range("D2").Value = "2016-01-28 04:28:13"
$ws3.Range("E2").Value = "a7e92d36-29de-4fd2-b571-1cfa984205c8.md"
$ws3.Range("F2").Value = "a7e92d36-29de-4fd2-b571-1cfa984205c8.aa6171980a82ab00837cc07c12fda249f142ae56.de-de.xlf"
$ws3.Range("G2").Value = "2016-01-28 04:28:59"

$ws3.Range("A3").Value = "ffffff4d4b21dd-2fbd-4145-b0e7-ea00cb7bc0b8.md"

$ws3.Range("A4").Value = "08ecc5f1-0974-4d29-aea9-f8af491d9dda.md"
$ws3.Range("B4").Value = "Ready for handoff"
$ws3.Range("C4").Value = "08ecc5f1-0974-4d29-aea9-f8af491d9dda.4f1445b64f0f841ffe719951dfaafccd1dce7119.de-de.xlf"
$ws3.Range("D4").Value = "2016-01-28 04:33:43"
$ws3.Range("E4").Value = "08ecc5f1-0974-4d29-aea9-f8af491d9dda.md"
$ws3.Range("F4").Value = "08ecc5f1-0974-4d29-aea9-f8af491d9dda.4f1445b64f0f841ffe719951dfaafccd1dce7119.de-de.xlf"
$ws3.Range("G4").Value = "2016-01-28 04:32:54"
